$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds plain-text values (e.g. "1.000", "26.619.45") that must
# stay text rather than being auto-coerced to numbers by Excel, so force the
# column to a text format before writing any of the new price strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.603.76'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '1.854.65'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '264.52'
$ws.Range("E5").Value = '  +2.38%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.5230'
$ws.Range("E7").Value = '  -0.65%  '
$ws.Range("D8").Value = '0.3282'
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").Value = '0.06800'
$ws.Range("D10").Value = '18.84'
$ws.Range("E10").Value = '  -3.13%  '
$ws.Range("D11").Value = '0.7772'
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = '0.07765'
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").Value = '1.843.22'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("D14").Value = '88.59'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("D15").Value = '5.026'
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '14.02'
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '0.000007975'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '0.9997'
$ws.Range("D20").Value = '26.633.88'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '2.083.92'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").Value = '4.641'
$ws.Range("E22").Value = '  +0.77%  '
$ws.Range("D23").Value = '9.556'
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").Value = '6.000'
$ws.Range("D25").Value = '144.40'
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("D26").Value = '2.197'
$ws.Range("E26").Value = '  -6.62%  '
$ws.Range("D27").Value = '1.677'
$ws.Range("E27").Value = '  +2.29%  '
$ws.Range("D28").Value = '17.02'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '112.36'
$ws.Range("E29").Value = '  +1.09%  '
$ws.Range("D30").Value = '4.201'
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("D31").Value = '4.157'
$ws.Range("E31").Value = '  -0.96%  '
$ws.Range("D32").Value = '0.08756'
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("D33").Value = '0.04832'
$ws.Range("E33").Value = '  -0.51%  '
$ws.Range("D34").Value = '1.138'
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.867'
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7147'
$ws.Range("E36").Value = '  +0.75%  '
$ws.Range("D37").Value = '3.118'
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = '0.01787'
$ws.Range("D39").Value = '2.205'
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("D40").Value = '0.4883'
$ws.Range("E40").Value = '  -1.64%  '
$ws.Range("D41").Value = '112.92'
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("D42").Value = '0.9012'
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").Value = '6.091'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = '0.9996'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '7.728'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").Value = '0.4197'
$ws.Range("E46").Value = '  -2.32%  '
$ws.Range("D47").Value = '9.150'
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").Value = '0.05930'
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").Value = '0.1245'
$ws.Range("E49").Value = '  -3.73%  '
$ws.Range("D50").Value = '35.01'
$ws.Range("E50").Value = '  -0.98%  '
$ws.Range("D51").Value = '0.8852'
$ws.Range("E51").Value = '  +2.32%  '
